# Update the "satimage 4435" results sheet (active sheet) with the second
# experiment's re-run cross-validation numbers for the H column, reshuffle
# the per-row "difference from best" formulas to key off column J (the new
# best-performing method) instead of column H, and refresh the view so the
# freshly highlighted B37:J66 block is what shows on reopen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated raw scores for column H (rows 2-31) -----------------------
$ws.Range("H2").Value  = 91.431792559200005
$ws.Range("H3").Value  = 91.995490417100001
$ws.Range("H4").Value  = 92.446448703499996
$ws.Range("H5").Value  = 91.882750845499999
$ws.Range("H6").Value  = 92.446448703499996
$ws.Range("H7").Value  = 91.882750845499999
$ws.Range("H8").Value  = 91.882750845499999
$ws.Range("H9").Value  = 92.333709131899994
$ws.Range("H10").Value = 91.544532130799993
$ws.Range("H11").Value = 91.544532130799993
$ws.Range("H12").Value = 93.122886132999994
$ws.Range("H13").Value = 92.220969560300006
$ws.Range("H14").Value = 91.431792559200005
$ws.Range("H15").Value = 91.882750845499999
$ws.Range("H16").Value = 90.755355129700007
$ws.Range("H17").Value = 90.755355129700007
$ws.Range("H18").Value = 91.544532130799993
$ws.Range("H19").Value = 91.882750845499999
$ws.Range("H20").Value = 92.897406989900006
$ws.Range("H21").Value = 92.220969560300006
$ws.Range("H22").Value = 91.431792559200005
$ws.Range("H23").Value = 91.657271702399996
$ws.Range("H24").Value = 92.784667418300003
$ws.Range("H25").Value = 90.868094701199993
$ws.Range("H26").Value = 91.995490417100001
$ws.Range("H27").Value = 91.882750845499999
$ws.Range("H28").Value = 90.642615558100005
$ws.Range("H29").Value = 91.995490417100001
$ws.Range("H30").Value = 91.995490417100001
$ws.Range("H31").Value = 92.559188275099999

# --- Row 37: per-column "gap vs. best" now anchored on column J --------
$ws.Range("B37").Formula = "=J2-B2"
$ws.Range("C37").Formula = "=J2-C2"
$ws.Range("D37").Formula = "=J2-D2"
$ws.Range("E37").Formula = "=J2-E2"
$ws.Range("F37").Formula = "=J2-F2"
$ws.Range("G37").Formula = "=J2-G2"
$ws.Range("H37").Formula = "=J2-H2"
$ws.Range("I37").Formula = "=J2-I2"
$ws.Range("J37").Formula = "=J2-J2"

# --- Rows 38:66: same shift, filled down as shared formulas ------------
$ws.Range("B38:B66").FormulaR1C1 = "=R[-35]C[8]-R[-35]C[0]"
$ws.Range("C38:C66").FormulaR1C1 = "=R[-35]C[7]-R[-35]C[0]"
$ws.Range("D38:D66").FormulaR1C1 = "=R[-35]C[6]-R[-35]C[0]"
$ws.Range("E38:E66").FormulaR1C1 = "=R[-35]C[5]-R[-35]C[0]"
$ws.Range("F38:F66").FormulaR1C1 = "=R[-35]C[4]-R[-35]C[0]"
$ws.Range("G38:G66").FormulaR1C1 = "=R[-35]C[3]-R[-35]C[0]"
$ws.Range("H38:H66").FormulaR1C1 = "=R[-35]C[2]-R[-35]C[0]"
$ws.Range("I38:I66").FormulaR1C1 = "=R[-35]C[1]-R[-35]C[0]"
$ws.Range("J38:J66").FormulaR1C1 = "=R[-35]C[0]-R[-35]C[0]"

# --- View: zoom to 85% and leave the new summary block selected ---------
$excel.ActiveWindow.Zoom = 85
$ws.Range("B37:J66").Select()
